# Update automatico via Actualizar 04-14-2021 20-33-03
#
# The "Fecha" (date) column D holds one timestamp per batch of 14 rows
# (one row per monitored service/link). Every refresh cycle pushes the
# availability-check history down one slot:
#   - rows 2-15  (newest) get the timestamp of the run that just finished
#   - rows 16-29 get what used to be rows 2-15's timestamp
#   - rows 30-43 get what used to be rows 16-29's timestamp
# (rows 30-43's previous timestamp ages out and is discarded).
#
# The serial date values below are the exact OLE-Automation-date doubles
# produced by that refresh, taken from the source-of-truth update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$batchSize     = 14
$firstDataRow  = 2

# Newest-to-oldest: each entry is the serial date to stamp across its
# 14-row batch, starting at row 2.
$batchValues = @(
    44300.85556826066,   # rows 2-15  (this run's timestamp, brand new)
    44267.74495982639,   # rows 16-29 (was rows 2-15's timestamp)
    44267.72340784722    # rows 30-43 (was rows 16-29's timestamp)
)

for ($b = 0; $b -lt $batchValues.Length; $b++) {
    $startRow = $firstDataRow + ($b * $batchSize)
    $value = $batchValues[$b]

    for ($i = 0; $i -lt $batchSize; $i++) {
        $row = $startRow + $i
        $ws.Cells.Item($row, 4).Value2 = $value
    }
}
